$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1690
$ws.Range("I2").Value = 1275
$ws.Range("K2").Value = 1275
$ws.Range("M2").Value = -1162
$ws.Range("H12").Value = 480.5
$ws.Range("J12").Value = 505
$ws.Range("L12").Value = 505
$ws.Range("N12").Value = -845
$ws.Range("H43").Value = 2966.3333
$ws.Range("I43").Value = 2966.3333
$ws.Range("K43").Value = 2966.3333
$ws.Range("M43").Value = -2897.3333
$ws.Range("H86").Value = 10464.571
$ws.Range("I86").Value = 4692.7144
$ws.Range("J86").Value = 16236.429
$ws.Range("K86").Value = 4692.7144
$ws.Range("L86").Value = 16236.429
$ws.Range("M86").Value = -3569.7144
$ws.Range("N86").Value = -18482.429
$ws.Range("H89").Value = 10464.571
$ws.Range("I89").Value = 4692.7144
$ws.Range("J89").Value = 16236.429
$ws.Range("K89").Value = 23463.572
$ws.Range("L89").Value = 81182.145
$ws.Range("M89").Value = -17847.572
$ws.Range("N89").Value = -92414.145
$ws.Range("H100").Value = 2027.5714
$ws.Range("I100").Value = 2678.8
$ws.Range("K100").Value = 2678.8
$ws.Range("M100").Value = -2137.8
$ws.Range("H107").Value = 1343.6364
$ws.Range("I107").Value = 722.04
$ws.Range("J107").Value = 3286.125
$ws.Range("K107").Value = 722.04
$ws.Range("L107").Value = 3286.125
$ws.Range("M107").Value = 1197.96
$ws.Range("N107").Value = -7126.125
$ws.Range("H135").Value = 1516.3334
$ws.Range("I135").Value = 1595.3334
$ws.Range("J135").Value = 1476.8334
$ws.Range("K135").Value = 14358.0006
$ws.Range("L135").Value = 13291.5006
$ws.Range("M135").Value = -11823.0006
$ws.Range("N135").Value = -18361.5006
$ws.Range("H138").Value = 4114.846
$ws.Range("I138").Value = 3595
$ws.Range("J138").Value = 4158.1665
$ws.Range("K138").Value = 10785
$ws.Range("L138").Value = 12474.4995
$ws.Range("M138").Value = -5645
$ws.Range("N138").Value = -22754.4995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13289.343
$ws.Range("I32").Value = 10761.454
$ws.Range("K32").Value = 10761.454
$ws.Range("M32").Value = -10474.454
$ws.Range("H37").Value = 24500
$ws.Range("J37").Value = 24500
$ws.Range("L37").Value = 24500
$ws.Range("N37").Value = -25046
$ws.Range("H39").Value = 3250
$ws.Range("I39").Value = 3250
$ws.Range("K39").Value = 3250
$ws.Range("M39").Value = -2730
$ws.Range("H44").Value = 28750
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 35000
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -35976
$ws.Range("H63").Value = 2853
$ws.Range("I63").Value = 2105.1428
$ws.Range("K63").Value = 2105.1428
$ws.Range("M63").Value = -1419.1428
$ws.Range("H66").Value = 2853
$ws.Range("I66").Value = 2105.1428
$ws.Range("K66").Value = 10525.714
$ws.Range("M66").Value = -7093.714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 1000
$ws.Range("J24").Value = 1000
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = -1470
$ws.Range("H94").Value = 1199.4286
$ws.Range("I94").Value = 956.7368
$ws.Range("J94").Value = 3505
$ws.Range("K94").Value = 956.7368
$ws.Range("L94").Value = 3505
$ws.Range("M94").Value = -505.7368
$ws.Range("N94").Value = -4407
$ws.Range("H107").Value = 2407.6
$ws.Range("I107").Value = 1346
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 1346
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 574
$ws.Range("N107").Value = -7840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 594.7778
$ws.Range("I8").Value = 104.5
$ws.Range("J8").Value = 987
$ws.Range("K8").Value = 104.5
$ws.Range("L8").Value = 987
$ws.Range("M8").Value = 35.5
$ws.Range("N8").Value = -1267
$ws.Range("H41").Value = 17240.416
$ws.Range("J41").Value = 17988.5
$ws.Range("L41").Value = 17988.5
$ws.Range("N41").Value = -18844.5
$ws.Range("H59").Value = 29195.691
$ws.Range("J59").Value = 34998.57
$ws.Range("L59").Value = 34998.57
$ws.Range("N59").Value = -37288.57
$ws.Range("H99").Value = 7899.857
$ws.Range("J99").Value = 7499.8
$ws.Range("L99").Value = 7499.8
$ws.Range("N99").Value = -10495.8
$ws.Range("H107").Value = 777.625
$ws.Range("I107").Value = 503
$ws.Range("K107").Value = 503
$ws.Range("M107").Value = 1417
$ws.Range("H122").Value = 2663.1428
$ws.Range("I122").Value = 2663.1428
$ws.Range("K122").Value = 7989.428400000001
$ws.Range("M122").Value = -5539.428400000001
$ws.Range("H126").Value = 7899.857
$ws.Range("J126").Value = 7499.8
$ws.Range("L126").Value = 22499.4
$ws.Range("N126").Value = -27439.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 179.6
$ws.Range("I2").Value = 84.90909000000001
$ws.Range("K2").Value = 509.4545400000001
$ws.Range("M2").Value = -396.4545400000001
$ws.Range("H12").Value = 113
$ws.Range("I12").Value = 70.5
$ws.Range("J12").Value = 141.33333
$ws.Range("K12").Value = 211.5
$ws.Range("L12").Value = 423.99999
$ws.Range("M12").Value = -38.5
$ws.Range("N12").Value = -769.99999
$ws.Range("H137").Value = 5561.5557
$ws.Range("I137").Value = 1380.75
$ws.Range("K137").Value = 4142.25
$ws.Range("M137").Value = 957.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 23325.285
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 32355.4
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 32355.4
$ws.Range("M43").Value = -599
$ws.Range("N43").Value = -32657.4
$ws.Range("H113").Value = 1085.5
$ws.Range("I113").Value = 1085.5
$ws.Range("K113").Value = 1085.5
$ws.Range("M113").Value = 1084.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3134.1538
$ws.Range("H132").Value = 10562.45
$ws.Range("I132").Value = 9013.053
$ws.Range("J132").Value = 40001
$ws.Range("K132").Value = 27039.159
$ws.Range("L132").Value = 120003
$ws.Range("M132").Value = -24509.159
$ws.Range("N132").Value = -125063
$ws.Range("H136").Value = 3252
$ws.Range("I136").Value = 3252
$ws.Range("K136").Value = 9756
$ws.Range("M136").Value = -7206

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4247.125
$ws.Range("I136").Value = 4247.125
$ws.Range("K136").Value = 12741.375
$ws.Range("M136").Value = -10191.375
